$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5437
$ws1.Range("F5").Value = 62
$ws1.Range("F6").Value = 76

# Sheet "全部类型" (All Types)
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F3").Value = 5437
$ws2.Range("F6").Value = 62
$ws2.Range("F7").Value = 76
